$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40"; change it to the text "1".
# Format the cell as Text first so the numeric-looking value "1" is kept
# as a string (matching a new shared-string entry) instead of being
# auto-converted to a number.
$r = $ws.Range("B11")
$r.NumberFormat = "@"
$r.Value = "1"
